# FBW A380X Checklist: fix typo "FKEX" -> "FLEX" in the
# "V1/VR/V2/FLEX TEMP. . . . . . . . . . . . . . CHECK" checklist line.
$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "V1/VR/V2/FKEX TEMP",  # FindText
    $true,                 # MatchCase
    $false,                # MatchWholeWord
    $false,                # MatchWildcards
    $false,                # MatchSoundsLike
    $false,                # MatchAllWordForms
    $true,                 # Forward
    1,                     # Wrap (wdFindContinue)
    $false,                # Format
    "V1/VR/V2/FLEX TEMP",  # ReplaceWith
    2                      # Replace (wdReplaceAll)
)

if (-not $found) {
    throw "Could not locate the 'V1/VR/V2/FKEX TEMP' checklist text to fix."
}
